$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "40.175.53"
$ws.Range("E2").Value = "  +3.49%  "
$ws.Range("D3").Value = "2.241.53"
$ws.Range("E3").Value = "  +1.53%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "296.33"
$ws.Range("E5").Value = "  +0.30%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "87.48"
$ws.Range("E6").Value = "  +10.25%  "
$ws.Range("E7").Value = "  +2.47%  "
$ws.Range("E8").Value = "  -0.08%  "
$ws.Range("E9").Value = "  +4.65%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "31.74"
$ws.Range("E10").Value = "  +14.51%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0801"
$ws.Range("E11").Value = "  +3.94%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "47.31"
$ws.Range("E12").Value = "  +3.62%  "
$ws.Range("E13").Value = "  +1.25%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.52"
$ws.Range("E14").Value = "  +7.46%  "
$ws.Range("D15").Value = "2.587.40"
$ws.Range("E15").Value = "  +1.66%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.29"
$ws.Range("E16").Value = "  +3.16%  "
$ws.Range("D17").Value = "2.237.30"
$ws.Range("E17").Value = "  -0.23%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.733"
$ws.Range("E18").Value = "  +3.56%  "
$ws.Range("D19").Value = "40.105.31"
$ws.Range("E19").Value = "  +3.55%  "
$ws.Range("E20").Value = "  +4.78%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.85"
$ws.Range("E21").Value = "  +2.81%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.93"
$ws.Range("E22").Value = "  +12.43%  "
$ws.Range("E23").Value = "  +1.71%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "236.37"
$ws.Range("E24").Value = "  +5.22%  "
$ws.Range("E25").Value = "  -0.02%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.48"
$ws.Range("E26").Value = "  +4.32%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.86"
$ws.Range("E27").Value = "  +7.18%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "23.02"
$ws.Range("E28").Value = "  +4.17%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.23"
$ws.Range("E29").Value = "  +2.48%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.28"
$ws.Range("E30").Value = "  +4.42%  "
$ws.Range("E31").Value = "  +9.43%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "152.78"
$ws.Range("E32").Value = "  +3.75%  "
$ws.Range("E33").Value = "  +0.07%  "
$ws.Range("E34").Value = "  +4.14%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0724"
$ws.Range("E35").Value = "  +5.86%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.38"
$ws.Range("E36").Value = "  +3.07%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "16.65"
$ws.Range("E37").Value = "  +16.09%  "
$ws.Range("B38").Value = "Stellar"
$ws.Range("C38").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.112"
$ws.Range("E38").Value = "  +4.10%  "
$ws.Range("B39").Value = "Kaspa"
$ws.Range("C39").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.101"
$ws.Range("E39").Value = "  +6.97%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.73"
$ws.Range("E40").Value = "  +3.57%  "
$ws.Range("E41").Value = "  +8.79%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.88"
$ws.Range("E42").Value = "  +8.86%  "
$ws.Range("D43").Value = "2.042.60"
$ws.Range("E43").Value = "  +7.30%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.25"
$ws.Range("E44").Value = "  +10.09%  "
$ws.Range("E45").Value = "  +8.13%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.00"
$ws.Range("E46").Value = "  +13.09%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "16.37"
$ws.Range("E47").Value = "  +1.91%  "
$ws.Range("E48").Value = "  +3.44%  "
$ws.Range("D49").Value = "2.458.65"
$ws.Range("E49").Value = "  +1.52%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "71.55"
$ws.Range("E50").Value = "  +2.28%  "
$ws.Range("E51").Value = "  +15.90%  "
